$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the style of the existing header row (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Values for the new "Save" column (H2:H14), per row, taken from the diff
$saveValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
